$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update numeric and text cell values ---
$ws.Range("C5").Value = 59.51
$ws.Range("D5").Value = 54.94
$ws.Range("E5").Value = 54.94
$ws.Range("F5").Value = 41.14
$ws.Range("C6").Value = 1.37
$ws.Range("D6").Value = 2.49
$ws.Range("E6").Value = 2.49
$ws.Range("F6").Value = 2.49
$ws.Range("C7").Value = 1.96
$ws.Range("D7").Value = 1.43
$ws.Range("E7").Value = 1.43
$ws.Range("F7").Value = 1.43
$ws.Range("C8").Value = 37.16
$ws.Range("D8").Value = 41.14
$ws.Range("E8").Value = 41.14
$ws.Range("F8").Value = 54.94
$ws.Range("C17").Value = 3
$ws.Range("C18").Value = 4
$ws.Range("C28").Value = 0.33
$ws.Range("C29").Value = 0.67
$ws.Range("C30").Value = 0.44
$ws.Range("C31").Value = 0.56
$ws.Range("D38").Value = 62.2
$ws.Range("D39").Value = 42.7
$ws.Range("D40").Value = 34.55
$ws.Range("C41").Value = 0
$ws.Range("D41").Value = 42.07
$ws.Range("C42").Value = 14.12
$ws.Range("D42").Value = 43.26
$ws.Range("C50").Value = 0.14
$ws.Range("C51").Value = 0.19
$ws.Range("D51").Value = 0.02
$ws.Range("C60").Value = 0.1
$ws.Range("D60").Value = 1.39
$ws.Range("E60").Value = 8.36
$ws.Range("C61").Value = 3.25
$ws.Range("D61").Value = 3.67
$ws.Range("E61").Value = 102.93
$ws.Range("C62").Value = 4.1
$ws.Range("D62").Value = 6.5
$ws.Range("E62").Value = 73.31
$ws.Range("B72").Value = 4.1
$ws.Range("C72").Value = 6.5
$ws.Range("D72").Value = 73.31
$ws.Range("B73").Value = 20.7
$ws.Range("C73").Value = 21.19
$ws.Range("D73").Value = 113.54
$ws.Range("B74").Value = 8.55
$ws.Range("C74").Value = 16.2
$ws.Range("D74").Value = 61.34
$ws.Range("C82").Value = 2015
$ws.Range("D82").Value = 2014
$ws.Range("C83").Value = 6.44
$ws.Range("D83").Value = 6.44
$ws.Range("C84").Value = 0.64
$ws.Range("D84").Value = 0.64
$ws.Range("C85").Value = 0.47
$ws.Range("D85").Value = 0.47
$ws.Range("C96").Value = "FY 12/13"
$ws.Range("D96").Value = "FY 13/14"
$ws.Range("E96").Value = "FY 14/15"
$ws.Range("C97").Value = 0.3
$ws.Range("D97").Value = 6.44
$ws.Range("E97").Value = 6.29
$ws.Range("C98").Value = 0.05
$ws.Range("D98").Value = 0.64
$ws.Range("E98").Value = 0.84
$ws.Range("C99").Value = 0.47
$ws.Range("D99").Value = 0.47
$ws.Range("E99").Value = 0.47
$ws.Range("B110").Value = "FY -1/le"
$ws.Range("B111").Value = "FY -1/"
$ws.Range("B112").Value = "FY -1/"
$ws.Range("B113").Value = "FY -1/"
$ws.Range("B114").Value = "FY -1/"
$ws.Range("B115").Value = "FY -1/"
$ws.Range("D123").Value = 0
$ws.Range("D124").Value = 8.2
$ws.Range("D125").Value = 0
$ws.Range("C134").Value = "FY -1/"
$ws.Range("C135").Value = "FY -1/"
$ws.Range("C136").Value = "FY -1/"
$ws.Range("C137").Value = "FY -1/"
$ws.Range("C138").Value = "FY -1/"
$ws.Range("B484").Value = "''"
$ws.Range("B485").Value = "''"

# --- Clear cells that become empty ---
$ws.Range("C110").ClearContents()
$ws.Range("D110").ClearContents()
$ws.Range("E110").ClearContents()
$ws.Range("C111").ClearContents()
$ws.Range("D111").ClearContents()
$ws.Range("E111").ClearContents()
$ws.Range("C112").ClearContents()
$ws.Range("D112").ClearContents()
$ws.Range("E112").ClearContents()
$ws.Range("C113").ClearContents()
$ws.Range("D113").ClearContents()
$ws.Range("E113").ClearContents()
$ws.Range("C114").ClearContents()
$ws.Range("D114").ClearContents()
$ws.Range("E114").ClearContents()
$ws.Range("C115").ClearContents()
$ws.Range("D115").ClearContents()
$ws.Range("E115").ClearContents()
$ws.Range("D134").ClearContents()
$ws.Range("D135").ClearContents()
$ws.Range("D136").ClearContents()
$ws.Range("D137").ClearContents()
$ws.Range("D138").ClearContents()
$ws.Range("C153").ClearContents()
$ws.Range("D153").ClearContents()
$ws.Range("E153").ClearContents()
$ws.Range("C154").ClearContents()
$ws.Range("D154").ClearContents()
$ws.Range("E154").ClearContents()
$ws.Range("C155").ClearContents()
$ws.Range("D155").ClearContents()
$ws.Range("E155").ClearContents()
$ws.Range("C156").ClearContents()
$ws.Range("D156").ClearContents()
$ws.Range("E156").ClearContents()
$ws.Range("C157").ClearContents()
$ws.Range("D157").ClearContents()
$ws.Range("E157").ClearContents()
$ws.Range("C484").ClearContents()
$ws.Range("D484").ClearContents()
$ws.Range("E484").ClearContents()
$ws.Range("C485").ClearContents()
$ws.Range("D485").ClearContents()
$ws.Range("E485").ClearContents()
